# Updated after reviewed by Ruchi Ma'am
# Text/grammar corrections to Task Description cells (column D) in both
# story blocks, plus a planned-hours correction for the "Buffer Time"
# rows (E18/E35: 6 -> 2), which ripples through the SUM()/shared-formula
# cells automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Story block 1 (rows 3-18, SSDMS-11)
# ---------------------------------------------------------------------
$ws.Range("D3").Value  = "Understand the 'WHY' of the story."
$ws.Range("D4").Value  = "Understand forward and backward linkages."
$ws.Range("D6").Value  = "Finalze the required fields of TP,AB and grouping of fields in Accordions"
$ws.Range("D7").Value  = "Develop technical understanding"
$ws.Range("D11").Value = "Design page structure using HTML"
$ws.Range("D12").Value = "Add Design using Bootstrap"
$ws.Range("D13").Value = "Integrate and creating fields using AngularJS"
$ws.Range("D14").Value = "Add CSS to the page"
$ws.Range("D17").Value = "Incorporate Code Review changes"

# Buffer Time planned hours corrected from 6 to 2
$ws.Range("E18").Value = 2

# ---------------------------------------------------------------------
# Story block 2 (rows 21-35, SSDMS-51)
# ---------------------------------------------------------------------
$ws.Range("D21").Value = "Understand the 'WHY' of the story."
$ws.Range("D22").Value = "Understand forward and backward linkages."
$ws.Range("D26").Value = "ii) Collaborate with the frontend team(deciding the controller's name and type of data to be outputted to every query)"
$ws.Range("D27").Value = "iii) Write dynamic SQL queries "
$ws.Range("D28").Value = "Design the Controllers,Services for the SQL queries"
$ws.Range("D29").Value = "Create DTO's for all Queries"
$ws.Range("D30").Value = "Create DAO's for all Queries"
$ws.Range("D31").Value = "Debugg"
$ws.Range("D34").Value = "Incorporate Code Review changes"

# Buffer Time planned hours corrected from 6 to 2
$ws.Range("E35").Value = 2

# ---------------------------------------------------------------------
# Update the view/selection to match the reviewed state (active cell
# moved to D8, top area scrolled down a bit).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("D8").Select()
